$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("DK2").Value = 90
$ws.Range("DL2").Value = "T"
$ws.Range("DL3").Value = "HG"
$ws.Range("IY3").Value = 90
$ws.Range("IZ3").Value = "T"
$ws.Range("ME3").Value = 90
$ws.Range("MF3").Value = "T"
$ws.Range("DL4").Value = "NR"
$ws.Range("IU4").Value = 90
$ws.Range("IV4").Value = "T"
$ws.Range("DL5").Value = "HG"
$ws.Range("DL6").Value = "HG"
$ws.Range("DK7").Value = 90
$ws.Range("DL7").Value = "T"
$ws.Range("DL8").Value = "HG"
$ws.Range("IU8").Value = 90
$ws.Range("IV8").Value = "T"
$ws.Range("IY8").Value = 90
$ws.Range("IZ8").Value = "T"
$ws.Range("DK9").Value = 90
$ws.Range("DL9").Value = "T"
$ws.Range("DL10").Value = "HG"
$ws.Range("DK11").Value = 66
$ws.Range("DL11").Value = "T"
$ws.Range("DL12").Value = "HG"
$ws.Range("DL13").Value = "HG"
$ws.Range("DK14").Value = 77
$ws.Range("DL14").Value = "T"
$ws.Range("DK15").Value = 77
$ws.Range("DL15").Value = "T"
$ws.Range("DK16").Value = 90
$ws.Range("DL16").Value = "T"
$ws.Range("DK17").Value = 13
$ws.Range("DL17").Value = "R"
$ws.Range("DK18").Value = 6
$ws.Range("DL18").Value = "R"
$ws.Range("DL19").Value = "HG"
$ws.Range("IY19").Value = 90
$ws.Range("IZ19").Value = "T"
$ws.Range("DK20").Value = 84
$ws.Range("DL20").Value = "T"
$ws.Range("DL21").Value = "HG"
$ws.Range("IU21").Value = 68
$ws.Range("IV21").Value = "T"
$ws.Range("DK22").Value = 90
$ws.Range("DL22").Value = "T"
$ws.Range("DL23").Value = "HG"
$ws.Range("IU23").Value = 90
$ws.Range("IV23").Value = "T"
$ws.Range("ME23").Value = 90
$ws.Range("MF23").Value = "T"
$ws.Range("DK24").Value = 90
$ws.Range("DL24").Value = "T"
$ws.Range("DL25").Value = "HG"
$ws.Range("IY25").Value = 90
$ws.Range("IZ25").Value = "T"
$ws.Range("DK26").Value = 13
$ws.Range("DL26").Value = "T"
$ws.Range("DL27").Value = "HG"
$ws.Range("DK28").Value = 24
$ws.Range("DL28").Value = "T"
$ws.Range("IU28").Value = 90
$ws.Range("IV28").Value = "T"
$ws.Range("DK29").Value = 90
$ws.Range("DL29").Value = "T"

$ws.Range("MJ20").Select()
